$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (shared strings "carrier" -> "carrier2", "volume" -> "volume2")
$ws.Range("A1").Value = "carrier2"
$ws.Range("B1").Value = "volume2"

# Fix audiogram: rows 25-31 column A should be 700 instead of 70
for ($r = 25; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = 700
}

# Update current selection to F10
$ws.Range("F10").Select()
